# The workbook gained a new data row: a row was inserted at sheet row 141,
# pushing the former rows 141-270 down to 142-271, and the new row 141 was
# populated with a fresh observation (Espinaca, Femacal de La Calera).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 141, shifting rows 141:270 down to 142:271.
$ws.Rows.Item(141).Insert()

# Populate the newly inserted row 141 with the new observation's data.
$ws.Cells.Item(141, 1).Value = 3
$ws.Cells.Item(141, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(141, 3).Value = "Coquimbo"
$ws.Cells.Item(141, 4).Value = 44586
$ws.Cells.Item(141, 5).Value = 5
$ws.Cells.Item(141, 6).Value = 100112012
$ws.Cells.Item(141, 7).Value = "Espinaca"
$ws.Cells.Item(141, 8).Value = "Sin especificar"
$ws.Cells.Item(141, 9).Value = "Primera"
$ws.Cells.Item(141, 10).Value = 140
$ws.Cells.Item(141, 11).Value = 4000
$ws.Cells.Item(141, 12).Value = 4300
$ws.Cells.Item(141, 13).Value = 4150
$ws.Cells.Item(141, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(141, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(141, 16).Value = 1383
$ws.Cells.Item(141, 17).Value = 3
$ws.Cells.Item(141, 18).Value = "Hortaliza"
